# GUI for editing member information
#
# Updates the "Attendance" sheet:
#  - drops the old "Paid" column (C) header/values, replacing it with a
#    one-off PAID note for the member who just paid (row 3)
#  - fixes a couple of placeholder member numbers/names
#  - appends two new attendance-date columns (X = August 08 2016,
#    Y = August 09 2016) with each member's check-in time for those dates

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Attendance")

# --- remove the old "Paid" column's data (header + its two values) ---
$ws.Range("C1").ClearContents()
$ws.Range("C2").ClearContents()
$ws.Range("C4").ClearContents()

# --- member data corrections ---
# Row 3: fix typo'd name, and this member has now paid
$ws.Range("B3").Value = "java script"
$ws.Range("C3").Value = "PAID"

# Row 4: member number correction + name tweak
$ws.Range("A4").Style = "Normal"
$ws.Range("A4").Value = 543543543
$ws.Range("A4").NumberFormat = "@"
$ws.Range("B4").Value = "darkest coffees"

# Row 7: replace placeholder number/name with the real member data
$ws.Range("A7").Style = "Normal"
$ws.Range("A7").Value = 604604604
$ws.Range("A7").NumberFormat = "@"
$ws.Range("B7").Value = "d code"

# Row 8: replace placeholder number/name with the real member data
$ws.Range("A8").Style = "Normal"
$ws.Range("A8").Value = 789789789
$ws.Range("A8").NumberFormat = "@"
$ws.Range("B8").Value = "turtwig overgrow"

# --- new attendance columns for August 08 2016 (X) and August 09 2016 (Y) ---
$ws.Range("X1").Value = "August 08 2016"
$ws.Range("Y1").Value = "August 09 2016"

$ws.Range("X2").Value = "03:56 PM"
$ws.Range("Y2").Value = "02:44 PM"

$ws.Range("X3").Value = "03:57 PM"
$ws.Range("Y3").Value = "02:45 PM"

$ws.Range("X6").Value = "03:58 PM"

$ws.Range("Y7").Value = "02:44 PM"

$ws.Range("Y8").Value = "02:44 PM"

# --- cosmetic: widen the Name column a bit, size the new date column like
#     the other attendance-time columns ---
$ws.Columns.Item(2).ColumnWidth = 14.75
$ws.Columns.Item(23).ColumnWidth = 13.25

# --- view state: move the selection, zoom out a bit ---
$ws.Activate()
$ws.Range("F10").Select()
$excel.ActiveWindow.Zoom = 48

Write-Output "edit complete"
